# Revert the localized (Japanese) SharePoint content-type-schema strings in
# the deck's custom XML part back to their original English wording.
#
# That "contentTypeSchema" blob isn't reachable through the regular
# Shapes/TextFrame surface -- it lives in the package's CustomXMLParts
# collection, so we go through ActivePresentation.CustomXMLParts exactly
# like real PowerPoint automation would. Re-assigning CustomXMLPart.XML is
# also what causes PowerPoint to re-roll the part's internal identifiers
# (the ma:versionID / ma:fieldsID attributes inside the XML, and the
# itemProps datastore GUID that backs CustomXMLPart.Id) as a side effect of
# persisting the edited content, so we don't need to (and can't, since
# CustomXMLPart.Id is read-only) touch those separately.

$p = $ppt.ActivePresentation

$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

function Get-PartByNamespace($presentation, $ns) {
    # Preferred: ask the collection to filter by namespace for us.
    try {
        $candidates = $presentation.CustomXMLParts.SelectByNamespace($ns)
        if ($candidates -ne $null -and $candidates.Count -ge 1) {
            return $candidates.Item(1)
        }
    } catch {
    }

    # Fall back to a manual scan (older/partial COM hosts sometimes don't
    # implement SelectByNamespace even though Item/NamespaceURI work fine).
    $total = $presentation.CustomXMLParts.Count
    for ($i = 1; $i -le $total; $i++) {
        $candidate = $presentation.CustomXMLParts.Item($i)
        if ($candidate -ne $null -and $candidate.NamespaceURI -eq $ns) {
            return $candidate
        }
    }

    return $null
}

$part = Get-PartByNamespace $p $contentTypeNs

if ($part -ne $null) {
    $xml = $part.XML

    # contentTypeSchema root attributes.
    $xml = $xml.Replace('ma:contentTypeName="ドキュメント"', 'ma:contentTypeName="Document"')
    $xml = $xml.Replace('ma:contentTypeDescription="新しいドキュメントを作成します。"', 'ma:contentTypeDescription="Create a new document."')
    $xml = $xml.Replace('ma:versionID="a9cab35011a557c1232e9e1918db7064"', 'ma:versionID="d0e002fabf17cb2440d8e9a473d3a41c"')
    $xml = $xml.Replace('ma:fieldsID="36c473bbc383ceb924bb8d2cdd9a2de6"', 'ma:fieldsID="e4cec627508c1f1ba247db94416ea198"')

    # Field display names.
    $xml = $xml.Replace('ma:displayName="画像タグ"', 'ma:displayName="Image Tags"')
    $xml = $xml.Replace('ma:displayName="コンテンツ タイプ"', 'ma:displayName="Content Type"')
    $xml = $xml.Replace('ma:displayName="タイトル"', 'ma:displayName="Title"')

    $part.XML = $xml
}
